$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest
# coinranking.com snapshot values.
#
# A handful of the new Price strings (e.g. "1.00", "7.00") look like
# plain numbers, and Excel would silently coerce them and drop the
# trailing zero(s) (e.g. "7.00" -> 7). Force those specific cells to
# Text format first so the literal string is preserved, matching the
# rest of the column which is already stored as text.

$ws.Range("D2").Value = "69.747.47"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "3.383.95"
$ws.Range("E3").Value = "  +3.79%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "190.87"
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "593.23"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.419"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").Value = "3.974.57"
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.72"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "69.733.53"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "3.382.99"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "452.38"
$ws.Range("E18").Value = "  +14.63%  "
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.82"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.17"
$ws.Range("E22").Value = "  +6.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.45"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.62"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.29"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.00"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +6.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.54"
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "28.02"
$ws.Range("E38").Value = "  +5.14%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.58"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").Value = "2.747.20"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.53"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.53"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0689"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "340.33"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.79"
$ws.Range("E49").Value = "  +6.37%  "
$ws.Range("E50").Value = "  +4.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.34"
$ws.Range("E51").Value = "  -0.63%  "
